$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 8 (CANSECO LEAL ANGELA): clean up the garbled Tutor name in column H.
$ws.Range("H8").Value = "KARINA JAQUELINE LEAL RIVERA"

# 2) Row 11 (CRUZ LOPEZ AISHA NAOMI): fill in the previously-missing contact
#    details (Correo, Tel_Movil, Tel_Fijo, Tutor, Telefono_Tutor). The phone
#    number columns must stay text, so they are entered with a leading
#    apostrophe (forces text-typed entry, same as a user typing it in Excel)
#    and then ClearFormats() strips the resulting "quote prefix" cell style
#    back off so the cell is left with no style override, matching its
#    untouched siblings.
$ws.Range("E11").Value = "aisha.naomi05@hotmail.com"

$ws.Range("F11").Value = "'2721538846"
$ws.Range("F11").ClearFormats()

$ws.Range("G11").Value = "'2721538846"
$ws.Range("G11").ClearFormats()

$ws.Range("H11").Value = "JACINTO CRUZ MARTÍNEZ"

$ws.Range("J11").Value = "'2721757221"
$ws.Range("J11").ClearFormats()

# 3) Row 20 (MATA CANSECO CRISTIAN ARTURO): the Correo_Tutor cell held the
#    bogus value "Kakakabs" - clear it out entirely.
$ws.Range("I20").ClearContents()
